$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.444.81"
$ws.Range("E2").Value = "  -3.02%  "
$ws.Range("D3").Value = "2.976.74"
$ws.Range("E3").Value = "  -4.93%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "496.01"
$ws.Range("D6").Value = "134.63"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.972.26"
$ws.Range("E8").Value = "  -5.08%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.426"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "7.18"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -7.29%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.127"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.490.53"
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "25.05"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "56.454.14"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.981.73"
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000146"
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "5.78"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.76"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "324.40"
$ws.Range("E22").Value = "  -5.63%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.467"
$ws.Range("E24").Value = "  -7.71%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "61.48"
$ws.Range("E25").Value = "  -9.48%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0893"
$ws.Range("E28").Value = "  -5.76%  "
$ws.Range("B29").Value = "USDe"
$ws.Range("C29").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "6.53"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.74"
$ws.Range("E33").Value = "  -6.19%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "20.19"
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "153.07"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.45"
$ws.Range("E36").Value = "  -7.25%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  -6.27%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "5.59"
$ws.Range("E38").Value = "  -10.21%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0669"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "23.14"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.013.78"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "36.57"
$ws.Range("E42").Value = "  -9.58%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.638"
$ws.Range("E44").Value = "  -7.41%  "
$ws.Range("D45").Value = "0.990"
$ws.Range("E45").Value = "  -8.73%  "
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.215.97"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "3.56"
$ws.Range("E48").Value = "  -8.55%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0238"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "5.75"
$ws.Range("E51").Value = "  -6.86%  "

Write-Output "Applied cryptos update"
